# Insert a new weekly price record as row 64, shifting existing rows 64-68 down to 65-69.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 44776
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112009
$ws.Cells.Item(64, 7).Value = "Acelga"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 200
$ws.Cells.Item(64, 11).Value = 1800
$ws.Cells.Item(64, 12).Value = 2000
$ws.Cells.Item(64, 13).Value = 1900
$ws.Cells.Item(64, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 633
$ws.Cells.Item(64, 17).Value = 3
$ws.Cells.Item(64, 18).Value = "Hortaliza"
